$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correct an existing value: G161 (Pm2.5) 34 -> 35 ---
$ws.Cells.Item(161, 7).Value = 35

# --- Append the new daily readings as rows 162-167 ---
# Columns: A=data(calendarisitca) date, B=data(zile) day index, C=Co, D=ica, E=No2, F=Pm10, G=Pm2.5, H=So2
$newData = @(
    @(44047, 161, 362, 111, 2, 2, 29, 6),
    @(44048, 162, 365, 106, 3, 3, 30, 7),
    @(44049, 163, 374, 110, 4, 4, 31, 9),
    @(44050, 164, 388, 107, 5, 5, 29, 12),
    @(44051, 165, 365, 115, 6, 6, 28, 15),
    @(44052, 166, 419, 109, 3, 3, 43, 8)
)

$startRow = 162
for ($i = 0; $i -lt $newData.Length; $i++) {
    $r = $startRow + $i
    $row = $newData[$i]

    # Match the formatting used on the row directly above (date style on A, centered style on B)
    $ws.Cells.Item(161, 1).Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)
    $ws.Cells.Item(161, 2).Copy()
    $ws.Cells.Item($r, 2).PasteSpecial(-4122)

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
}

# --- Column L picked up an explicit width (incidental, no data lives there) ---
$ws.Columns.Item(12).ColumnWidth = 8

# --- Scroll the view down so the newly entered rows are visible ---
$excel.ActiveWindow.ScrollRow = 151
